# Update Leve profit-tracking data (currentAveragePrice* / LeveProfit* columns)
# sourced from the scheduled market-data refresh, per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3339.2144
$ws.Range("I40").Value = 2122.6365
$ws.Range("K40").Value = 2122.6365
$ws.Range("M40").Value = -1947.6365
$ws.Range("H51").Value = 1499.5
$ws.Range("I51").Value = 1495
$ws.Range("K51").Value = 1495
$ws.Range("M51").Value = -1011
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("H106").Value = 4166.5
$ws.Range("I106").Value = 1333
$ws.Range("K106").Value = 1333
$ws.Range("M106").Value = -702
$ws.Range("H127").Value = 3253.5
$ws.Range("I127").Value = 398.5
$ws.Range("K127").Value = 1195.5
$ws.Range("M127").Value = 3764.5
$ws.Range("H129").Value = 1249.75
$ws.Range("I129").Value = 1249.75
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 3749.25
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 1250.75
$ws.Range("N129").ClearContents()
$ws.Range("H137").Value = 1768
$ws.Range("I137").Value = 1631.5834
$ws.Range("K137").Value = 4894.7502
$ws.Range("M137").Value = -2344.7502

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2009.3
$ws.Range("I45").Value = 1974.25
$ws.Range("J45").Value = 2149.5
$ws.Range("K45").Value = 1974.25
$ws.Range("L45").Value = 2149.5
$ws.Range("M45").Value = -1597.25
$ws.Range("N45").Value = -2903.5
$ws.Range("H132").Value = 3221.2144
$ws.Range("I132").Value = 2344.5557
$ws.Range("K132").Value = 7033.6671
$ws.Range("M132").Value = -4503.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1286.4445
$ws.Range("I107").Value = 1058.6923
$ws.Range("K107").Value = 1058.6923
$ws.Range("M107").Value = 861.3077000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 596.3333
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 989
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 989
$ws.Range("M22").Value = -50
$ws.Range("N22").Value = -1689
$ws.Range("H31").Value = 2042
$ws.Range("I31").Value = 1681.5
$ws.Range("K31").Value = 1681.5
$ws.Range("M31").Value = -1386.5
$ws.Range("H34").Value = 2042
$ws.Range("I34").Value = 1681.5
$ws.Range("K34").Value = 1681.5
$ws.Range("M34").Value = -1479.5
$ws.Range("H62").Value = 6749.8335
$ws.Range("I62").Value = 4166.6665
$ws.Range("J62").Value = 9333
$ws.Range("K62").Value = 4166.6665
$ws.Range("L62").Value = 9333
$ws.Range("M62").Value = -3542.6665
$ws.Range("N62").Value = -10581
$ws.Range("H65").Value = 6749.8335
$ws.Range("I65").Value = 4166.6665
$ws.Range("J65").Value = 9333
$ws.Range("K65").Value = 20833.3325
$ws.Range("L65").Value = 46665
$ws.Range("M65").Value = -17713.3325
$ws.Range("N65").Value = -52905
$ws.Range("H105").Value = 4136.5
$ws.Range("I105").Value = 2704.75
$ws.Range("J105").Value = 7000
$ws.Range("K105").Value = 2704.75
$ws.Range("L105").Value = 7000
$ws.Range("M105").Value = -957.75
$ws.Range("N105").Value = -10494
$ws.Range("H122").Value = 5650.4287
$ws.Range("I122").Value = 5152.8
$ws.Range("K122").Value = 15458.4
$ws.Range("M122").Value = -13008.4
$ws.Range("H132").Value = 4812.125
$ws.Range("I132").Value = 4633
$ws.Range("K132").Value = 13899
$ws.Range("M132").Value = -11369
$ws.Range("H134").Value = 2840.5
$ws.Range("I134").Value = 2450
$ws.Range("K134").Value = 7350
$ws.Range("M134").Value = -4815

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6057.15
$ws.Range("I3").Value = 2361.353
$ws.Range("J3").Value = 27000
$ws.Range("K3").Value = 7084.059
$ws.Range("L3").Value = 81000
$ws.Range("M3").Value = -6972.059
$ws.Range("N3").Value = -81224
$ws.Range("H86").Value = 1284.75
$ws.Range("I86").Value = 2222
$ws.Range("J86").Value = 347.5
$ws.Range("K86").Value = 6666
$ws.Range("L86").Value = 1042.5
$ws.Range("M86").Value = -5480
$ws.Range("N86").Value = -3414.5
$ws.Range("H89").Value = 1284.75
$ws.Range("I89").Value = 2222
$ws.Range("J89").Value = 347.5
$ws.Range("K89").Value = 19998
$ws.Range("L89").Value = 3127.5
$ws.Range("M89").Value = -14070
$ws.Range("N89").Value = -14983.5
$ws.Range("H122").Value = 2899.1428
$ws.Range("J122").Value = 3291.25
$ws.Range("L122").Value = 29621.25
$ws.Range("N122").Value = -34521.25
$ws.Range("H129").Value = 1943.2727
$ws.Range("I129").Value = 1232.5
$ws.Range("J129").Value = 2349.4285
$ws.Range("K129").Value = 3697.5
$ws.Range("L129").Value = 7048.2855
$ws.Range("M129").Value = 1302.5
$ws.Range("N129").Value = -17048.2855
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1848.2
$ws.Range("I102").Value = 2060.25
$ws.Range("K102").Value = 2060.25
$ws.Range("M102").Value = -438.25
$ws.Range("H122").Value = 3667.6667
$ws.Range("I122").Value = 3251.5
$ws.Range("K122").Value = 9754.5
$ws.Range("M122").Value = -7304.5
$ws.Range("H132").Value = 4259
$ws.Range("I132").Value = 3769.5
$ws.Range("K132").Value = 11308.5
$ws.Range("M132").Value = -8778.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2074.25
$ws.Range("I68").Value = 2160.111
$ws.Range("J68").Value = 1816.6666
$ws.Range("K68").Value = 2160.111
$ws.Range("L68").Value = 1816.6666
$ws.Range("M68").Value = -1411.111
$ws.Range("N68").Value = -3314.6666
$ws.Range("H71").Value = 2074.25
$ws.Range("I71").Value = 2160.111
$ws.Range("J71").Value = 1816.6666
$ws.Range("K71").Value = 10800.555
$ws.Range("L71").Value = 9083.333000000001
$ws.Range("M71").Value = -7056.555
$ws.Range("N71").Value = -16571.333
$ws.Range("H93").Value = 1872.25
$ws.Range("I93").Value = 1639.1333
$ws.Range("J93").Value = 2260.7778
$ws.Range("K93").Value = 1639.1333
$ws.Range("L93").Value = 2260.7778
$ws.Range("M93").Value = -391.1333
$ws.Range("N93").Value = -4756.7778
$ws.Range("H132").Value = 1821.9166
$ws.Range("I132").Value = 986.5
$ws.Range("K132").Value = 2959.5
$ws.Range("M132").Value = -429.5
$ws.Range("H136").Value = 3830.4
$ws.Range("I136").Value = 3830.4
$ws.Range("K136").Value = 11491.2
$ws.Range("M136").Value = -8941.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1641.9474

Write-Host "Updated 176 cells across 8 sheets"